$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: new appointment reusing existing Manolache Alexandra / Dragos Truta / Tuns simplu values
$ws.Range("A7").Value = "2021-05-24 13:00"
$ws.Range("B7").Value = "Manolache Alexandra"
$ws.Range("C7").Value = "Dragos Truta"
$ws.Range("D7").Value = "Tuns simplu"
$ws.Range("E7").Value = 45.0

# Row 8: new appointment with a new employee (manager) and new service
$ws.Range("A8").Value = "2021-05-25 12:30"
$ws.Range("B8").Value = "Giuredea Manuela-Ioana"
$ws.Range("C8").Value = "Dragos Truta"
$ws.Range("D8").Value = "Pedichiura cu oja semipermanenta"
$ws.Range("E8").Value = 50.0
